$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 691331
$ws.Range("J96").Value = 1611449.9
$ws.Range("L96").Value = 4834349.699999999
$ws.Range("N96").Value = -4837095.699999999
$ws.Range("H98").Value = 4465852
$ws.Range("I98").Value = 4809055.5
$ws.Range("K98").Value = 4809055.5
$ws.Range("M98").Value = -4807557.5
$ws.Range("H122").Value = 4465852
$ws.Range("I122").Value = 4809055.5
$ws.Range("K122").Value = 14427166.5
$ws.Range("M122").Value = -14424716.5
$ws.Range("H137").Value = 2000.8918
$ws.Range("I137").Value = 1891.8928
$ws.Range("J137").Value = 2340
$ws.Range("K137").Value = 5675.678400000001
$ws.Range("L137").Value = 7020
$ws.Range("M137").Value = -3125.678400000001
$ws.Range("N137").Value = -12120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12979084
$ws.Range("I61").Value = 14331531
$ws.Range("K61").Value = 14331531
$ws.Range("M61").Value = -14331319
$ws.Range("H74").Value = 1565.275
$ws.Range("I74").Value = 1550.3823
$ws.Range("K74").Value = 1550.3823
$ws.Range("M74").Value = -676.3823
$ws.Range("H77").Value = 1565.275
$ws.Range("I77").Value = 1550.3823
$ws.Range("K77").Value = 7751.9115
$ws.Range("M77").Value = -3383.9115
$ws.Range("H97").Value = 1740.8125
$ws.Range("I97").Value = 1373.3077
$ws.Range("K97").Value = 1373.3077
$ws.Range("M97").Value = -877.3077000000001
$ws.Range("H136").Value = 12979084
$ws.Range("I136").Value = 14331531
$ws.Range("K136").Value = 42994593
$ws.Range("M136").Value = -42992043

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 511134.66
$ws.Range("I105").Value = 916773.0600000001
$ws.Range("J105").Value = 4086.7
$ws.Range("K105").Value = 916773.0600000001
$ws.Range("L105").Value = 4086.7
$ws.Range("M105").Value = -915026.0600000001
$ws.Range("N105").Value = -7580.7
$ws.Range("H134").Value = 9093485
$ws.Range("I134").Value = 2833.7
$ws.Range("K134").Value = 8501.099999999999
$ws.Range("M134").Value = -5966.099999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19610392
$ws.Range("I31").Value = 35717150
$ws.Range("J31").Value = 2162
$ws.Range("K31").Value = 35717150
$ws.Range("L31").Value = 2162
$ws.Range("M31").Value = -35716855
$ws.Range("N31").Value = -2752
$ws.Range("H34").Value = 19610392
$ws.Range("I34").Value = 35717150
$ws.Range("J34").Value = 2162
$ws.Range("K34").Value = 35717150
$ws.Range("L34").Value = 2162
$ws.Range("M34").Value = -35716948
$ws.Range("N34").Value = -2566
$ws.Range("H58").Value = 3792
$ws.Range("I58").Value = 3536.2
$ws.Range("J58").Value = 4005.1667
$ws.Range("K58").Value = 3536.2
$ws.Range("L58").Value = 4005.1667
$ws.Range("M58").Value = -3333.2
$ws.Range("N58").Value = -4411.1667
$ws.Range("H105").Value = 2643.4546
$ws.Range("I105").Value = 2133
$ws.Range("K105").Value = 2133
$ws.Range("M105").Value = -386
$ws.Range("H107").Value = 2668.9167
$ws.Range("J107").Value = 3578
$ws.Range("L107").Value = 3578
$ws.Range("N107").Value = -7418
$ws.Range("H132").Value = 3435.1765
$ws.Range("I132").Value = 3387.9285
$ws.Range("K132").Value = 10163.7855
$ws.Range("M132").Value = -7633.7855
$ws.Range("H136").Value = 3792
$ws.Range("I136").Value = 3536.2
$ws.Range("J136").Value = 4005.1667
$ws.Range("K136").Value = 10608.6
$ws.Range("L136").Value = 12015.5001
$ws.Range("M136").Value = -8058.599999999999
$ws.Range("N136").Value = -17115.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350.85
$ws.Range("I5").Value = 1359.8
$ws.Range("J5").Value = 1341.9
$ws.Range("K5").Value = 4079.4
$ws.Range("L5").Value = 4025.7
$ws.Range("M5").Value = -3967.4
$ws.Range("N5").Value = -4249.700000000001
$ws.Range("H38").Value = 994.5263
$ws.Range("J38").Value = 994.0714
$ws.Range("L38").Value = 2982.2142
$ws.Range("N38").Value = -3676.2142
$ws.Range("H135").Value = 1350.85
$ws.Range("I135").Value = 1359.8
$ws.Range("J135").Value = 1341.9
$ws.Range("K135").Value = 12238.2
$ws.Range("L135").Value = 12077.1
$ws.Range("M135").Value = -9703.199999999999
$ws.Range("N135").Value = -17147.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 7913.5
$ws.Range("I23").Value = 77.5
$ws.Range("J23").Value = 15749.5
$ws.Range("K23").Value = 77.5
$ws.Range("L23").Value = 15749.5
$ws.Range("M23").Value = 145.5
$ws.Range("N23").Value = -16195.5
$ws.Range("H40").Value = 59999.5
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29849
$ws.Range("H74").Value = 99999
$ws.Range("J74").Value = 99999
$ws.Range("L74").Value = 99999
$ws.Range("N74").Value = -101871
$ws.Range("H77").Value = 99999
$ws.Range("J77").Value = 99999
$ws.Range("L77").Value = 299997
$ws.Range("N77").Value = -309357
$ws.Range("H107").Value = 909.5599999999999
$ws.Range("J107").Value = 841.1429000000001
$ws.Range("L107").Value = 841.1429000000001
$ws.Range("N107").Value = -4681.1429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 9800
$ws.Range("I25").Value = 9800
$ws.Range("K25").Value = 9800
$ws.Range("M25").Value = -9570
$ws.Range("H40").Value = 3413
$ws.Range("I40").Value = 3290.923
$ws.Range("K40").Value = 3290.923
$ws.Range("M40").Value = -3154.923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 257218.5
$ws.Range("I3").Value = 4435
$ws.Range("K3").Value = 4435
$ws.Range("M3").Value = -4321
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -5480
$ws.Range("H113").Value = 598.7
$ws.Range("I113").Value = 514
$ws.Range("J113").Value = 725.75
$ws.Range("K113").Value = 1542
$ws.Range("L113").Value = 2177.25
$ws.Range("M113").Value = 628
$ws.Range("N113").Value = -6517.25
$ws.Range("H126").Value = 2710.125
$ws.Range("I126").Value = 2781.8845
$ws.Range("K126").Value = 8345.6535
$ws.Range("M126").Value = -5875.6535

